# Apply updated dSF (column F) values as part of a data repull / mean recalculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    9  = -1
    12 = -2
    15 = 0
    16 = 1
    17 = -1
    35 = 1
    49 = -4
    51 = 4
    54 = 1
    57 = 4
    58 = -2
    67 = 0
    71 = -1
    81 = 0
    83 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
